# Weekly update: insert a new data row for the most recent week at row 414,
# pushing all subsequent rows (previously 414-502) down by one (415-503).
# This mirrors how the source "Consolidado" feed prepends the latest week's
# observation to the historical series on each refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 414; existing rows 414:502 shift to 415:503
# and the sheet's dimension grows from A1:T502 to A1:T503 automatically.
$ws.Rows("414:414").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(414, 1).Value  = 4
$ws.Cells.Item(414, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(414, 3).Value  = "Los Lagos"
$ws.Cells.Item(414, 4).Value  = 44711
$ws.Cells.Item(414, 5).Value  = 10
$ws.Cells.Item(414, 6).Value  = "Fruta"
$ws.Cells.Item(414, 7).Value  = 100108
$ws.Cells.Item(414, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(414, 9).Value  = 100108006
$ws.Cells.Item(414, 10).Value = "Plátano"
$ws.Cells.Item(414, 11).Value = "Sin especificar"
$ws.Cells.Item(414, 12).Value = "Primera Pintón"
$ws.Cells.Item(414, 13).Value = 500
$ws.Cells.Item(414, 14).Value = 18000
$ws.Cells.Item(414, 15).Value = 18000
$ws.Cells.Item(414, 16).Value = 18000
$ws.Cells.Item(414, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(414, 18).Value = "Ecuador"
$ws.Cells.Item(414, 19).Value = 900
$ws.Cells.Item(414, 20).Value = 20
